# Add "2022-Q4" data:
#  1. Duplicate the "2022-Q3" sheet, place the copy right before it, rename
#     to "2022-Q4", and update its fund-size / position figures.
#  2. Insert a corresponding summary row at the top of the "总计" sheet and
#     renumber the existing index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q4" worksheet from a copy of "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3Index = $q3.Index
$q3.Copy($q3, $null)             # new copy is placed immediately before $q3,
                                  # taking over $q3's former tab position
$q4 = $wb.Worksheets.Item($q3Index)
$q4.Name = "2022-Q4"

# Row 2 (005585 / 银河文体娱乐主题灵活配置混合A)
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "3.15"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "88.90"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "4.86"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.1531"

# Row 3 (015667 / 银河文体娱乐主题灵活配置混合C)
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "0.38"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "88.90"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "4.86"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0185"

# ---------------------------------------------------------------------
# 2. Update the "总计" overview sheet with the new quarter's totals
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Remember the existing rows 2-4 before they get shifted down.
$rowsData = @()
for ($r = 2; $r -le 4; $r++) {
    $rowsData += , @(
        $total.Cells.Item($r, 2).Value2,
        $total.Cells.Item($r, 3).Value2,
        $total.Cells.Item($r, 4).Value2
    )
}

# The index column's formatting (bold, centred, boxed) needs to reach the
# brand-new row 5 too - seed it now by copying row 4's format down.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)   # xlPasteFormats

# New row 2: 2022-Q4 summary.
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 0.17

# Shift the previously-read rows down by one, renumbering column A.
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = 3 + $i
    $total.Cells.Item($r, 1).Value = $i + 1
    $total.Cells.Item($r, 2).Value = $rowsData[$i][0]
    $total.Cells.Item($r, 3).Value = $rowsData[$i][1]
    $total.Cells.Item($r, 4).Value = $rowsData[$i][2]
}

# Copying "2022-Q3" made the new sheet the active tab; restore the
# original active sheet ("2020-Q4", the last tab) as before the edit.
$wb.Worksheets.Item("2020-Q4").Activate()
